$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = '34.472.45'
$ws.Cells.Item(2, 5).Value = '  -0.07%  '

# Row 3
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = '1.804.63'
$ws.Cells.Item(3, 5).Value = '  -0.03%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  -0.12%  '

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '224.84'
$ws.Cells.Item(5, 5).Value = '  -1.31%  '

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '0.605'

# Row 7
$ws.Cells.Item(7, 5).Value = '  -0.04%  '

# Row 8
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '39.11'
$ws.Cells.Item(8, 5).Value = '  +6.38%  '

# Row 9
$ws.Cells.Item(9, 5).Value = '  -3.31%  '

# Row 10
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '0.0669'
$ws.Cells.Item(10, 5).Value = '  -3.39%  '

# Row 11
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '0.0983'
$ws.Cells.Item(11, 5).Value = '  +2.08%  '

# Row 12
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '2.068.74'
$ws.Cells.Item(12, 5).Value = '  +0.17%  '

# Row 13
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '11.04'
$ws.Cells.Item(13, 5).Value = '  -4.62%  '

# Row 14
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '1.798.21'
$ws.Cells.Item(14, 5).Value = '  -0.47%  '

# Row 15
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '0.630'
$ws.Cells.Item(15, 5).Value = '  -2.48%  '

# Row 16
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '34.486.89'
$ws.Cells.Item(16, 5).Value = '  +0.12%  '

# Row 17
$ws.Cells.Item(17, 5).Value = '  -2.18%  '

# Row 18
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '68.21'
$ws.Cells.Item(18, 5).Value = '  -3.03%  '

# Row 19
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '239.78'
$ws.Cells.Item(19, 5).Value = '  -2.35%  '

# Row 20
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '0.0₃0769'
$ws.Cells.Item(20, 5).Value = '  -2.79%  '

# Row 21
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '11.14'
$ws.Cells.Item(21, 5).Value = '  -4.02%  '

# Row 23
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '4.09'
$ws.Cells.Item(23, 5).Value = '  -2.57%  '

# Row 24
$ws.Cells.Item(24, 5).Value = '  +0.57%  '

# Row 25
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '170.85'
$ws.Cells.Item(25, 5).Value = '  -0.84%  '

# Row 26
$ws.Cells.Item(26, 2).Value = 'Cosmos'
$ws.Cells.Item(26, 3).Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '7.70'
$ws.Cells.Item(26, 5).Value = '  -3.87%  '

# Row 27
$ws.Cells.Item(27, 2).Value = 'EthereumClassic'
$ws.Cells.Item(27, 3).Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '17.62'
$ws.Cells.Item(27, 5).Value = '  +3.91%  '

# Row 28
$ws.Cells.Item(28, 5).Value = '  +3.66%  '

# Row 29
$ws.Cells.Item(29, 5).Value = '  -0.07%  '

# Row 30
$ws.Cells.Item(30, 5).Value = '  -1.31%  '

# Row 31
$ws.Cells.Item(31, 5).Value = '  -2.37%  '

# Row 32
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '0.0514'
$ws.Cells.Item(32, 5).Value = '  -2.87%  '

# Row 33
$ws.Cells.Item(33, 5).Value = '  -4.32%  '

# Row 34
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '1.81'
$ws.Cells.Item(34, 5).Value = '  +0.13%  '

# Row 35
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '0.640'
$ws.Cells.Item(35, 5).Value = '  -5.18%  '

# Row 36
$ws.Cells.Item(36, 2).Value = 'TrustWalletToken'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '1.06'
$ws.Cells.Item(36, 5).Value = '  -0.75%  '

# Row 37
$ws.Cells.Item(37, 2).Value = 'Maker'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '1.304.28'
$ws.Cells.Item(37, 5).Value = '  -6.32%  '

# Row 38
$ws.Cells.Item(38, 5).Value = '  -2.02%  '

# Row 39
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '2.30'
$ws.Cells.Item(39, 5).Value = '  -6.32%  '

# Row 40
$ws.Cells.Item(40, 5).Value = '  +1.99%  '

# Row 41
$ws.Cells.Item(41, 5).Value = '  +0.68%  '

# Row 42
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '82.17'
$ws.Cells.Item(42, 5).Value = '  -1.08%  '

# Row 43
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '0.951'
$ws.Cells.Item(43, 5).Value = '  -1.15%  '

# Row 44
$ws.Cells.Item(44, 5).Value = '  -1.33%  '

# Row 45
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '13.91'
$ws.Cells.Item(45, 5).Value = '  +3.19%  '

# Row 46
$ws.Cells.Item(46, 5).Value = '  +2.40%  '

# Row 47
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '1.968.96'
$ws.Cells.Item(47, 5).Value = '  +0.12%  '

# Row 48
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '5.77'
$ws.Cells.Item(48, 5).Value = '  -4.17%  '

# Row 49
$ws.Cells.Item(49, 5).Value = '  -0.11%  '

# Row 50
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '102.43'
$ws.Cells.Item(50, 5).Value = '  -1.98%  '

# Row 51
$ws.Cells.Item(51, 2).Value = 'BabyDogeCoin'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '0.0₆0120'
$ws.Cells.Item(51, 5).Value = '  -6.04%  '
